# Update co-occurrence theme pairs / frequencies per the lemmatization refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 3).Value = "Environnement"
$ws.Cells.Item(2, 4).Value = 11

# Row 3
$ws.Cells.Item(3, 2).Value = "Emploi"
$ws.Cells.Item(3, 3).Value = "Justice"
$ws.Cells.Item(3, 4).Value = 11

# Row 4
$ws.Cells.Item(4, 4).Value = 10

# Row 5
$ws.Cells.Item(5, 3).Value = "Économie"
$ws.Cells.Item(5, 4).Value = 9

# Row 6
$ws.Cells.Item(6, 2).Value = "Gouvernance"
$ws.Cells.Item(6, 4).Value = 8

# Row 7
$ws.Cells.Item(7, 2).Value = "Emploi"
$ws.Cells.Item(7, 4).Value = 8

# Row 8
$ws.Cells.Item(8, 2).Value = "Environnement"
$ws.Cells.Item(8, 3).Value = "Gouvernance"
$ws.Cells.Item(8, 4).Value = 8

# Row 9
$ws.Cells.Item(9, 2).Value = "Emploi"
$ws.Cells.Item(9, 3).Value = "Santé"
$ws.Cells.Item(9, 4).Value = 8

# Row 10
$ws.Cells.Item(10, 2).Value = "Environnement"
$ws.Cells.Item(10, 4).Value = 8

# Row 11
$ws.Cells.Item(11, 2).Value = "Gouvernance"
$ws.Cells.Item(11, 3).Value = "Économie"
$ws.Cells.Item(11, 4).Value = 8

# Row 12
$ws.Cells.Item(12, 4).Value = 51

# Row 13
$ws.Cells.Item(13, 2).Value = "Gouvernance"
$ws.Cells.Item(13, 4).Value = 40

# Row 14
$ws.Cells.Item(14, 2).Value = "Justice"
$ws.Cells.Item(14, 3).Value = "Économie"
$ws.Cells.Item(14, 4).Value = 33

# Row 15
$ws.Cells.Item(15, 2).Value = "Environnement"
$ws.Cells.Item(15, 3).Value = "Gouvernance"
$ws.Cells.Item(15, 4).Value = 32

# Row 16
$ws.Cells.Item(16, 2).Value = "Social"
$ws.Cells.Item(16, 3).Value = "Économie"
$ws.Cells.Item(16, 4).Value = 30

# Row 17
$ws.Cells.Item(17, 2).Value = "Environnement"
$ws.Cells.Item(17, 3).Value = "Justice"
$ws.Cells.Item(17, 4).Value = 29

# Row 18
$ws.Cells.Item(18, 2).Value = "Gouvernance"
$ws.Cells.Item(18, 3).Value = "Social"
$ws.Cells.Item(18, 4).Value = 27

# Row 19
$ws.Cells.Item(19, 2).Value = "Justice"
$ws.Cells.Item(19, 4).Value = 25

# Row 20
$ws.Cells.Item(20, 2).Value = "Gouvernance"
$ws.Cells.Item(20, 3).Value = "Justice"
$ws.Cells.Item(20, 4).Value = 24

# Row 21
$ws.Cells.Item(21, 2).Value = "Infrastructure"
$ws.Cells.Item(21, 3).Value = "Économie"
$ws.Cells.Item(21, 4).Value = 23

# Row 22
$ws.Cells.Item(22, 2).Value = "Environnement"
$ws.Cells.Item(22, 3).Value = "Gouvernance"
$ws.Cells.Item(22, 4).Value = 11

# Row 23
$ws.Cells.Item(23, 2).Value = "Gouvernance"
$ws.Cells.Item(23, 3).Value = "Justice"
$ws.Cells.Item(23, 4).Value = 10

# Row 24
$ws.Cells.Item(24, 3).Value = "Social"
$ws.Cells.Item(24, 4).Value = 9

# Row 25
$ws.Cells.Item(25, 2).Value = "Environnement"
$ws.Cells.Item(25, 3).Value = "Justice"
$ws.Cells.Item(25, 4).Value = 7

# Row 26
$ws.Cells.Item(26, 2).Value = "Justice"
$ws.Cells.Item(26, 3).Value = "Social"
$ws.Cells.Item(26, 4).Value = 6

# Row 27
$ws.Cells.Item(27, 2).Value = "Environnement"
$ws.Cells.Item(27, 3).Value = "Social"
$ws.Cells.Item(27, 4).Value = 5

# Row 28
$ws.Cells.Item(28, 4).Value = 4

# Row 29
$ws.Cells.Item(29, 2).Value = "Gouvernance"
$ws.Cells.Item(29, 3).Value = "Économie"

# Row 30
$ws.Cells.Item(30, 2).Value = "Justice"
$ws.Cells.Item(30, 3).Value = "Économie"
$ws.Cells.Item(30, 4).Value = 3

# Row 31
$ws.Cells.Item(31, 4).Value = 3

# Row 32
$ws.Cells.Item(32, 4).Value = 21

# Row 33
$ws.Cells.Item(33, 4).Value = 20

# Row 34
$ws.Cells.Item(34, 4).Value = 19

# Row 35
$ws.Cells.Item(35, 2).Value = "Social"

# Row 36
$ws.Cells.Item(36, 3).Value = "Économie"
$ws.Cells.Item(36, 4).Value = 18

# Row 37
$ws.Cells.Item(37, 2).Value = "Gouvernance"
$ws.Cells.Item(37, 3).Value = "Justice"

# Row 38
$ws.Cells.Item(38, 2).Value = "Justice"

# Row 39
$ws.Cells.Item(39, 2).Value = "Emploi"
$ws.Cells.Item(39, 3).Value = "Gouvernance"
$ws.Cells.Item(39, 4).Value = 17

# Row 40
$ws.Cells.Item(40, 2).Value = "Droits_Femme"
$ws.Cells.Item(40, 3).Value = "Social"

# Row 41
$ws.Cells.Item(41, 3).Value = "Social"
$ws.Cells.Item(41, 4).Value = 16
